$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$v = $ws.Range("Z1").Value
Write-Host "Type: $($v.GetType().Name)"
Write-Host "Val: $v"
